# Update cryptocurrency price/volume data to reflect the latest GitHub Actions scrape.
# Row 37/38 coin identities are swapped (Maker <-> MXToken) in addition to value updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.389.21'
$ws.Range("E2").Value = '  -0.08%  '

# Row 3
$ws.Range("D3").Value = '1.845.55'
$ws.Range("E3").Value = '  -0.25%  '

# Row 4
$ws.Range("D4").Value = '''0.9981'
$ws.Range("E4").Value = '  -0.26%  '

# Row 5
$ws.Range("D5").Value = '''240.67'
$ws.Range("E5").Value = '  +0.02%  '

# Row 6
$ws.Range("D6").Value = '''0.6312'
$ws.Range("E6").Value = '  +0.43%  '

# Row 7
$ws.Range("D7").Value = '''0.9995'
$ws.Range("E7").Value = '  -0.19%  '

# Row 8
$ws.Range("D8").Value = '''0.07488'
$ws.Range("E8").Value = '  -2.11%  '

# Row 9
$ws.Range("E9").Value = '  +0.00%  '

# Row 10
$ws.Range("D10").Value = '''24.41'
$ws.Range("E10").Value = '  -1.68%  '

# Row 11
$ws.Range("E11").Value = '  -0.12%  '

# Row 12
$ws.Range("D12").Value = '1.845.71'
$ws.Range("E12").Value = '  -2.25%  '

# Row 13
$ws.Range("E13").Value = '  -0.67%  '

# Row 14
$ws.Range("E14").Value = '  -0.42%  '

# Row 15
$ws.Range("D15").Value = '''0.00001020'
$ws.Range("E15").Value = '  -5.01%  '

# Row 16
$ws.Range("D16").Value = '''82.14'
$ws.Range("E16").Value = '  -1.48%  '

# Row 17
$ws.Range("D17").Value = '''6.143'
$ws.Range("E17").Value = '  -0.44%  '

# Row 18
$ws.Range("D18").Value = '29.423.39'
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$ws.Range("D19").Value = '''228.80'
$ws.Range("E19").Value = '  +0.24%  '

# Row 20
$ws.Range("E20").Value = '  -0.08%  '

# Row 21
$ws.Range("D21").Value = '''0.9989'
$ws.Range("E21").Value = '  -0.23%  '

# Row 22
$ws.Range("D22").Value = '''7.442'
$ws.Range("E22").Value = '  -0.22%  '

# Row 23
$ws.Range("D23").Value = '''0.9998'

# Row 24
$ws.Range("D24").Value = '''158.98'
$ws.Range("E24").Value = '  +0.69%  '

# Row 25
$ws.Range("D25").Value = '''0.1377'
$ws.Range("E25").Value = '  -0.43%  '

# Row 26
$ws.Range("D26").Value = '''8.423'
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$ws.Range("E27").Value = '  -0.92%  '

# Row 28
$ws.Range("D28").Value = '''0.06341'
$ws.Range("E28").Value = '  +13.32%  '

# Row 29
$ws.Range("D29").Value = '''1.382'
$ws.Range("E29").Value = '  -0.48%  '

# Row 30
$ws.Range("D30").Value = '''1.474'
$ws.Range("E30").Value = '  +0.62%  '

# Row 31
$ws.Range("D31").Value = '''4.093'

# Row 32
$ws.Range("D32").Value = '''4.056'
$ws.Range("E32").Value = '  -0.19%  '

# Row 33
$ws.Range("E33").Value = '  -1.15%  '

# Row 34
$ws.Range("D34").Value = '''1.141'
$ws.Range("E34").Value = '  -1.93%  '

# Row 35
$ws.Range("D35").Value = '''0.6970'
$ws.Range("E35").Value = '  +0.11%  '

# Row 36
$ws.Range("E36").Value = '  -0.33%  '

# Row 37
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '''2.837'
$ws.Range("E37").Value = '  +3.96%  '

# Row 38
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.255.22'
$ws.Range("E38").Value = '  +2.14%  '

# Row 39
$ws.Range("D39").Value = '''0.01822'
$ws.Range("E39").Value = '  +0.98%  '

# Row 40
$ws.Range("E40").Value = '  +1.64%  '

# Row 41
$ws.Range("D41").Value = '''0.9082'
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("D42").Value = '''0.9993'
$ws.Range("E42").Value = '  -0.17%  '

# Row 43
$ws.Range("D43").Value = '2.008.99'
$ws.Range("E43").Value = '  -15.47%  '

# Row 44
$ws.Range("D44").Value = '''101.37'
$ws.Range("E44").Value = '  -1.13%  '

# Row 45
$ws.Range("E45").Value = '  +0.49%  '

# Row 46
$ws.Range("D46").Value = '''0.1173'
$ws.Range("E46").Value = '  +2.11%  '

# Row 47
$ws.Range("D47").Value = '''7.047'
$ws.Range("E47").Value = '  -2.15%  '

# Row 48
$ws.Range("D48").Value = '''0.00000000116'
$ws.Range("E48").Value = '  +1.04%  '

# Row 49
$ws.Range("D49").Value = '''9.042'
$ws.Range("E49").Value = '  +0.43%  '

# Row 50
$ws.Range("D50").Value = '''1.687'
$ws.Range("E50").Value = '  +0.42%  '

# Row 51
$ws.Range("E51").Value = '  -2.13%  '
